$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "glycan"
$ws.Range("B1").Value = "binding_score"
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
$ws.Range("F1").Value = "flexibility"
$ws.Range("G1").Value = "has_multi_node_motifs"

# Row 2
$ws.Range("A2").Value = "Fuc(a1-2)Gal(b1-4)GlcNAc"
$ws.Range("B2").Value = 0.7597051643243998
$ws.Range("C2").Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'GlcNAc(b1-1)']"
$ws.Range("D2").Value = "['Fuc(a1-2)Gal(b1-4)GlcNAc']"
$ws.Range("E2").Value = 7.020136041242929
$ws.Range("F2").Value = 1.339857268468498
$ws.Range("G2").Value = $true

# Row 3
$ws.Range("A3").Value = "Fuc(a1-2)Gal(b1-4)[Fuc(a1-3)]GlcNAc"
$ws.Range("B3").Value = 0.2129727877185961
$ws.Range("C3").Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'GlcNAc(b1-1)']"
$ws.Range("D3").Value = "['Fuc(a1-2)Gal(b1-4)GlcNAc']"
$ws.Range("E3").Value = 5.991065472972884
$ws.Range("F3").Value = 0.6054263107241069
$ws.Range("G3").Value = $true

# Row 4
$ws.Range("A4").Value = "Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc"
$ws.Range("B4").Value = -0.3927742735235619
$ws.Range("C4").Value = "['Fuc(a1-2)', 'Gal(b1-4)', 'GlcNAc(b1-3)']"
$ws.Range("D4").Value = "['Fuc(a1-2)Gal(b1-4)GlcNAc']"
$ws.Range("E4").Value = 5.806341297683876
$ws.Range("F4").Value = 1.79358205293686
$ws.Range("G4").Value = $true

# Styling: header row + column A (rows 2-4) get bold font, thin border,
# horizontal-center / vertical-top alignment.
# NOTE: multi-area ("A1:G1,A2:A4") ranges only apply formatting to the
# first area in this runtime, so style the header row and the A2:A4
# block as two separate contiguous ranges instead.
$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$colARange = $ws.Range("A2:A4")
$colARange.Font.Bold = $true
$colARange.HorizontalAlignment = -4108
$colARange.VerticalAlignment = -4160
$colARange.Borders.LineStyle = 1

$ws.Range("A1:G4").Select()
